$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1951754385964912
$ws.Range("C2").Value = 0.5504385964912281
$ws.Range("J2").Value = 0.01535087719298246
$ws.Range("O2").Value = 0.001096491228070175
$ws.Range("P2").Value = 0.1458333333333333
$ws.Range("S2").Value = 0.09210526315789473
$ws.Range("B3").Value = 0.01691729323308271
$ws.Range("C3").Value = 0.04135338345864661
$ws.Range("J3").Value = 0.03007518796992481
$ws.Range("P3").Value = 0.7481203007518797
$ws.Range("S3").Value = 0.1635338345864662
$ws.Range("J4").Value = 0.03870967741935484
$ws.Range("P4").Value = 0.6903225806451613
$ws.Range("S4").Value = 0.2709677419354839
$ws.Range("J5").Value = 0.2857142857142857
$ws.Range("P5").Value = 0.1428571428571428
$ws.Range("S5").Value = 0.5714285714285714
$ws.Range("B6").Value = 0.08368794326241134
$ws.Range("D6").Value = 0.005673758865248227
$ws.Range("E6").Value = 0.002836879432624114
$ws.Range("F6").Value = 0.06666666666666667
$ws.Range("J6").Value = 0.2397163120567376
$ws.Range("O6").Value = 0.02127659574468085
$ws.Range("Q6").Value = 0.1588652482269504
$ws.Range("R6").Value = 0.06950354609929078
$ws.Range("S6").Value = 0.3517730496453901
$ws.Range("B7").Value = 0.0840064620355412
$ws.Range("D7").Value = 0.02746365105008077
$ws.Range("E7").Value = 0.001615508885298869
$ws.Range("F7").Value = 0.06946688206785137
$ws.Range("J7").Value = 0.1373182552504039
$ws.Range("O7").Value = 0.01938610662358643
$ws.Range("Q7").Value = 0.1663974151857835
$ws.Range("R7").Value = 0.06946688206785137
$ws.Range("S7").Value = 0.4248788368336026
$ws.Range("B8").Value = 0.1051454138702461
$ws.Range("D8").Value = 0.01864280387770321
$ws.Range("F8").Value = 0.06263982102908278
$ws.Range("J8").Value = 0.09619686800894854
$ws.Range("O8").Value = 0.01640566741237882
$ws.Range("Q8").Value = 0.2080536912751678
$ws.Range("R8").Value = 0.08650260999254288
$ws.Range("S8").Value = 0.4064131245339299
$ws.Range("B9").Value = 0.090625
$ws.Range("D9").Value = 0.021875
$ws.Range("F9").Value = 0.07968749999999999
$ws.Range("J9").Value = 0.109375
$ws.Range("O9").Value = 0.0234375
$ws.Range("Q9").Value = 0.165625
$ws.Range("R9").Value = 0.1
$ws.Range("S9").Value = 0.409375
$ws.Range("B10").Value = 0.1046099290780142
$ws.Range("D10").Value = 0.02608915906788247
$ws.Range("E10").Value = 0.001266464032421479
$ws.Range("F10").Value = 0.06433637284701114
$ws.Range("J10").Value = 0.1162613981762918
$ws.Range("O10").Value = 0.02077001013171226
$ws.Range("Q10").Value = 0.2084599797365755
$ws.Range("R10").Value = 0.0797872340425532
$ws.Range("S10").Value = 0.378419452887538
$ws.Range("G11").Value = 0.1440162271805274
$ws.Range("J11").Value = 0.09533468559837728
$ws.Range("K11").Value = 0.2038539553752536
$ws.Range("L11").Value = 0.5456389452332657
$ws.Range("S11").Value = 0.01115618661257606
$ws.Range("G12").Value = 0.6892361111111112
$ws.Range("J12").Value = 0.2256944444444444
$ws.Range("K12").Value = 0.01388888888888889
$ws.Range("L12").Value = 0.02604166666666667
$ws.Range("S12").Value = 0.04513888888888889
$ws.Range("F13").Value = 0.01379310344827586
$ws.Range("G13").Value = 0.6551724137931034
$ws.Range("J13").Value = 0.2620689655172414
$ws.Range("S13").Value = 0.06896551724137931
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.1666666666666667
$ws.Range("S14").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.02989130434782609
$ws.Range("H15").Value = 0.1331521739130435
$ws.Range("I15").Value = 0.07608695652173914
$ws.Range("J15").Value = 0.3668478260869565
$ws.Range("K15").Value = 0.05706521739130434
$ws.Range("M15").Value = 0.006793478260869565
$ws.Range("O15").Value = 0.08152173913043478
$ws.Range("S15").Value = 0.2486413043478261
$ws.Range("F16").Value = 0.02419354838709677
$ws.Range("H16").Value = 0.1629032258064516
$ws.Range("I16").Value = 0.1032258064516129
$ws.Range("J16").Value = 0.3854838709677419
$ws.Range("K16").Value = 0.1064516129032258
$ws.Range("M16").Value = 0.008064516129032258
$ws.Range("O16").Value = 0.05967741935483871
$ws.Range("S16").Value = 0.15
$ws.Range("F17").Value = 0.01981599433828733
$ws.Range("H17").Value = 0.1719745222929936
$ws.Range("I17").Value = 0.1012031139419674
$ws.Range("J17").Value = 0.4048124557678698
$ws.Range("K17").Value = 0.1019108280254777
$ws.Range("M17").Value = 0.02335456475583864
$ws.Range("N17").Value = 0.0007077140835102619
$ws.Range("O17").Value = 0.06086341118188252
$ws.Range("S17").Value = 0.1153573956121727
$ws.Range("F18").Value = 0.02051282051282051
$ws.Range("H18").Value = 0.1726495726495726
$ws.Range("I18").Value = 0.08717948717948718
$ws.Range("J18").Value = 0.4051282051282051
$ws.Range("K18").Value = 0.1042735042735043
$ws.Range("M18").Value = 0.02051282051282051
$ws.Range("O18").Value = 0.05641025641025641
$ws.Range("S18").Value = 0.1333333333333333
$ws.Range("F19").Value = 0.01708343649418173
$ws.Range("H19").Value = 0.2002971032433771
$ws.Range("I19").Value = 0.08195097796484278
$ws.Range("J19").Value = 0.3651894033176529
$ws.Range("K19").Value = 0.1121564743748453
$ws.Range("M19").Value = 0.02352067343401832
$ws.Range("N19").Value = 0.00198068828918049
$ws.Range("O19").Value = 0.07303788066353058
$ws.Range("S19").Value = 0.1247833622183709
